# Update cryptos list with refreshed Price/Volume(1h) figures (and, for a
# couple of rows whose ranking order flipped, the Coin/Link text too).
# A leading apostrophe is used for purely-numeric-looking Price strings so
# Excel keeps them as text (matching the original inline-string cells)
# instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.513.40'
$ws.Range("E2").Value = '  +4.39%  '
$ws.Range("D3").Value = '3.252.82'
$ws.Range("E3").Value = '  +3.70%  '
$ws.Range("D5").Value = '''578.61'
$ws.Range("E5").Value = '  +2.49%  '
$ws.Range("D6").Value = '''181.50'
$ws.Range("E6").Value = '  +6.68%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -4.04%  '
$ws.Range("D9").Value = '3.250.13'
$ws.Range("E9").Value = '  +3.66%  '
$ws.Range("E10").Value = '  +4.93%  '
$ws.Range("D11").Value = '''6.77'
$ws.Range("E11").Value = '  +3.50%  '
$ws.Range("E12").Value = '  +5.63%  '
$ws.Range("D13").Value = '3.822.75'
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("D15").Value = '''28.49'
$ws.Range("E15").Value = '  +5.57%  '
$ws.Range("D16").Value = '67.496.56'
$ws.Range("E16").Value = '  +4.51%  '
$ws.Range("E17").Value = '  +3.00%  '
$ws.Range("D18").Value = '3.252.48'
$ws.Range("E19").Value = '  +2.54%  '
$ws.Range("E20").Value = '  +5.39%  '
$ws.Range("D21").Value = '''377.04'
$ws.Range("E21").Value = '  +6.35%  '
$ws.Range("E22").Value = '  +5.45%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = '''71.42'
$ws.Range("E24").Value = '  +4.43%  '
$ws.Range("E25").Value = '  +2.66%  '
$ws.Range("D26").Value = '''0.0000119'
$ws.Range("E26").Value = '  +1.70%  '
$ws.Range("D27").Value = '''9.60'
$ws.Range("E27").Value = '  +0.69%  '
$ws.Range("E28").Value = '  +3.55%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("D30").Value = '''5.73'
$ws.Range("E30").Value = '  +7.02%  '
$ws.Range("E31").Value = '  +4.44%  '
$ws.Range("D32").Value = '''22.60'
$ws.Range("E32").Value = '  +3.21%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  +6.11%  '
$ws.Range("E35").Value = '  +4.02%  '
$ws.Range("D36").Value = '''163.38'
$ws.Range("E36").Value = '  +6.37%  '
$ws.Range("E37").Value = '  +4.03%  '
$ws.Range("D38").Value = '''0.854'
$ws.Range("E39").Value = '  +7.38%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '''6.79'
$ws.Range("E40").Value = '  +13.13%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '''26.81'
$ws.Range("E41").Value = '  +3.13%  '
$ws.Range("D42").Value = '''2.62'
$ws.Range("E42").Value = '  +3.33%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '''4.49'
$ws.Range("E43").Value = '  +8.04%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = '''361.35'
$ws.Range("E44").Value = '  +12.69%  '
$ws.Range("D45").Value = '2.747.92'
$ws.Range("E45").Value = '  +3.78%  '
$ws.Range("D46").Value = '''25.44'
$ws.Range("E46").Value = '  +5.43%  '
$ws.Range("D47").Value = '''40.59'
$ws.Range("E47").Value = '  +4.08%  '
$ws.Range("E48").Value = '  +3.19%  '
$ws.Range("E49").Value = '  +2.71%  '
$ws.Range("E50").Value = '  +7.41%  '
$ws.Range("E51").Value = '  +0.53%  '
